$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 13 (pushes former rows 13-23 down to 14-24).
# This correctly carries the existing custom row heights down with the
# content, since each height is tied to the row that holds the label.
$ws.Rows.Item(13).Insert()

# The new row 13 holds the "Docentes responsáveis:" answer (previously
# mis-placed as the "Objetivos:" answer). Write the values first...
$ws.Range("B13").Value = "2143261 - André Luis Ferraz"
$ws.Range("C13").Value = "2143261 - André Luis Ferraz"

# ...then fix up B13/C13's style: Insert() copied the row-above's format
# (column A's style) across the whole new row, but B/C need their normal
# per-column styles. Borrow the correct formats from the row below, which
# still has the regular B/C styling.
$ws.Range("B14:C14").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
# Row 13 never had an "A" label - drop the stray formatted cell Insert()
# left behind there.
$ws.Range("A13").Clear()

# Objetivos: replace the (misplaced) answer with the real objectives text.
$ws.Range("B10").Value = "Propiciar ao aluno conhecimentos básicos da Química Inorgânica envolvida em processos biológicos."
$ws.Range("C10").Value = "Propiciar ao aluno conhecimentos básicos da Química Inorgânica envolvida em processos biológicos."

# Programa resumido: (now row 14) replace placeholder "Semestral" with the
# actual short syllabus text (PT).
$ws.Range("B14").Value = "Estrutura molecular e ligação química; Orbitais moleculares e as moléculas de O2 e N2; Ácidos, bases e a correlação com os ligantes dos metais em solução; Complexos metálicos - teoria do campo cristalino; Sistemas biológicos de transporte; Transporte de O2 e transferência de elétrons em sistemas biológicos; Processos catalíticos - ácido/base e oxido-redução em metaloproteínas."
$ws.Range("C14").Value = "Estrutura molecular e ligação química; Orbitais moleculares e as moléculas de O2 e N2; Ácidos, bases e a correlação com os ligantes dos metais em solução; Complexos metálicos - teoria do campo cristalino; Sistemas biológicos de transporte; Transporte de O2 e transferência de elétrons em sistemas biológicos; Processos catalíticos - ácido/base e oxido-redução em metaloproteínas."

# Programa: (now row 16) replace placeholder date value with the actual
# full syllabus text (PT).
$ws.Range("B16").Value = "1. Estrutura molecular e ligação química: Teoria de ligação de valência, estrutura de compostos com C, N, O; Relação entre estrutura e propriedades fisico-químicas2. Orbitais moleculares e as moléculas de O2 e N2: Limitações da teoria de ligação de valência, reatividade diferenciada de O2 e N2, relevância do O2 em sistemas biológicos, espécies reativas de oxigênio3. Ácidos, bases e a correlação com os ligantes dos metais em solução: Afinidade das bases por metais de transição, equilíbrio químico em sistemas biológicos4. Complexos metálicos - teoria do campo cristalino: Teoria do campo cristalino e os compostos octaédricos e tetraédricos; íons de metais de transição em sistemas biológicos5. Sistemas biológicos de transporte: Transporte de O2 em mamíferos, transferência de elétrons dependente de metaloproteínas;6. Processos catalíticos - ácido/base e oxido-redução em metaloproteínas: Proteínas contendo íon Zn2+, peroxidases, oxidases."
$ws.Range("C16").Value = "1. Estrutura molecular e ligação química: Teoria de ligação de valência, estrutura de compostos com C, N, O; Relação entre estrutura e propriedades fisico-químicas2. Orbitais moleculares e as moléculas de O2 e N2: Limitações da teoria de ligação de valência, reatividade diferenciada de O2 e N2, relevância do O2 em sistemas biológicos, espécies reativas de oxigênio3. Ácidos, bases e a correlação com os ligantes dos metais em solução: Afinidade das bases por metais de transição, equilíbrio químico em sistemas biológicos4. Complexos metálicos - teoria do campo cristalino: Teoria do campo cristalino e os compostos octaédricos e tetraédricos; íons de metais de transição em sistemas biológicos5. Sistemas biológicos de transporte: Transporte de O2 em mamíferos, transferência de elétrons dependente de metaloproteínas;6. Processos catalíticos - ácido/base e oxido-redução em metaloproteínas: Proteínas contendo íon Zn2+, peroxidases, oxidases."

# Método / Critério / Norma de recuperação answers were each shifted up by
# one label (fixing a pre-existing off-by-one mistake in the sheet), and
# Bibliografia gets a brand-new answer.
$ws.Range("B19").Value = "A avaliação será feita por meio de provas escritas."
$ws.Range("C19").Value = "A avaliação será feita por meio de provas escritas."

$ws.Range("B20").Value = "A Nota final (NF) será calculada da seguinte maneira:NF = (P1 + 2*P2)/3Sendo que para P2 a matéria será cumulativa do semestre."
$ws.Range("C20").Value = "A Nota final (NF) será calculada da seguinte maneira:NF = (P1 + 2*P2)/3Sendo que para P2 a matéria será cumulativa do semestre."

$ws.Range("B21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$ws.Range("C21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"

$ws.Range("B22").Value = "1. Atkins e Jones, Princípios de Química, 5a edição, Bookman, 20112. Shiver e Atikins, Química Inorgânica, 4a edição, Bookman, 2008"
$ws.Range("C22").Value = "1. Atkins e Jones, Princípios de Química, 5a edição, Bookman, 20112. Shiver e Atikins, Química Inorgânica, 4a edição, Bookman, 2008"

# Column layout: column A's width/style now only applies to column 1 (it
# used to span columns 1-2 before column B got its own explicit column
# definition below).
$ws.Columns.Item(1).ColumnWidth = 30.7109375
